# Auto-generated cell updates derived from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.996.88'
$ws.Range("E2").Value = '  -3.16%  '
$ws.Range("D3").Value = '3.362.44'
$ws.Range("E3").Value = '  -2.69%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'566.47"
$ws.Range("E5").Value = '  -2.56%  '
$ws.Range("D6").Value = "'148.99"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'8.00"
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("D11").Value = "'0.414"
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("D12").Value = '3.937.57'
$ws.Range("E12").Value = '  -2.71%  '
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").Value = "'28.05"
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '3.354.52'
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = '  -1.41%  '
$ws.Range("D17").Value = '61.064.01'
$ws.Range("E17").Value = '  -3.14%  '
$ws.Range("D18").Value = "'6.34"
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").Value = "'14.28"
$ws.Range("E19").Value = '  -2.46%  '
$ws.Range("D20").Value = "'8.83"
$ws.Range("E20").Value = '  -3.95%  '
$ws.Range("D21").Value = "'375.84"
$ws.Range("E21").Value = '  -3.33%  '
$ws.Range("D22").Value = "'75.19"
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").Value = "'0.560"
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("D25").Value = '3.500.38'
$ws.Range("E26").Value = '  -6.35%  '
$ws.Range("D27").Value = "'0.176"
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("D29").Value = "'7.42"
$ws.Range("E29").Value = '  -3.57%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'7.73"
$ws.Range("E31").Value = '  -4.21%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = "'2.07"
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("D33").Value = "'22.89"
$ws.Range("E33").Value = '  -1.90%  '
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("D35").Value = "'5.38"
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").Value = "'170.66"
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("E37").Value = '  -6.29%  '
$ws.Range("D38").Value = "'6.79"
$ws.Range("E38").Value = '  -3.80%  '
$ws.Range("D39").Value = "'28.85"
$ws.Range("E39").Value = '  -10.15%  '
$ws.Range("D40").Value = '3.399.27'
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("D41").Value = "'0.0745"
$ws.Range("E41").Value = '  -4.55%  '
$ws.Range("D42").Value = "'42.33"
$ws.Range("E42").Value = '  -1.41%  '
$ws.Range("D43").Value = "'0.760"
$ws.Range("E43").Value = '  -4.22%  '
$ws.Range("D44").Value = "'4.29"
$ws.Range("E44").Value = '  -1.52%  '
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("E46").Value = '  -6.21%  '
$ws.Range("D47").Value = '2.488.07'
$ws.Range("E47").Value = '  -2.97%  '
$ws.Range("D48").Value = "'6.67"
$ws.Range("E48").Value = '  -3.14%  '
$ws.Range("D49").Value = "'22.53"
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").Value = "'0.0262"
$ws.Range("E51").Value = '  -2.20%  '
